$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 55556524
$ws.Range("I6").Value = 62501084
$ws.Range("J6").Value = 51
$ws.Range("K6").Value = 187503252
$ws.Range("L6").Value = 153
$ws.Range("M6").Value = -187503140
$ws.Range("N6").Value = -377

$ws.Range("H28").Value = 4653.4
$ws.Range("I28").Value = 4653.4
$ws.Range("K28").Value = 4653.4
$ws.Range("M28").Value = -4168.4

$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76248

$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231240

$ws.Range("H113").Value = 14326
$ws.Range("J113").Value = 5126.5
$ws.Range("L113").Value = 5126.5
$ws.Range("N113").Value = -11634.5

$ws.Range("H116").Value = 16977.6
$ws.Range("I116").Value = 19998.5
$ws.Range("K116").Value = 19998.5
$ws.Range("M116").Value = -16556.5

$ws.Range("H132").Value = 1729.8572
$ws.Range("I132").Value = 1256.4412
$ws.Range("J132").Value = 3741.875
$ws.Range("K132").Value = 3769.3236
$ws.Range("L132").Value = 11225.625
$ws.Range("M132").Value = -1239.3236
$ws.Range("N132").Value = -16285.625

$ws.Range("H137").Value = 7156.45
$ws.Range("I137").Value = 8876.5
$ws.Range("K137").Value = 26629.5
$ws.Range("M137").Value = -24079.5

$ws.Range("H138").Value = 3668.0852
$ws.Range("I138").Value = 3132.6667
$ws.Range("J138").Value = 3794.8948
$ws.Range("K138").Value = 9398.000100000001
$ws.Range("L138").Value = 11384.6844
$ws.Range("M138").Value = -4258.000100000001
$ws.Range("N138").Value = -21664.6844

$ws.Range("H141").Value = 3718
$ws.Range("I141").Value = 3530
$ws.Range("K141").Value = 10590
$ws.Range("M141").Value = -5410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2521.923
$ws.Range("I61").Value = 2499.682
$ws.Range("J61").Value = 2644.25
$ws.Range("K61").Value = 2499.682
$ws.Range("L61").Value = 2644.25
$ws.Range("M61").Value = -2287.682
$ws.Range("N61").Value = -3068.25

$ws.Range("H110").Value = 423.0909
$ws.Range("I110").Value = 423.0909
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 423.0909
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1621.9091
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 1795.3414
$ws.Range("I132").Value = 1419.0938
$ws.Range("J132").Value = 3133.111
$ws.Range("K132").Value = 4257.2814
$ws.Range("L132").Value = 9399.332999999999
$ws.Range("M132").Value = -1727.2814
$ws.Range("N132").Value = -14459.333

$ws.Range("H136").Value = 2521.923
$ws.Range("I136").Value = 2499.682
$ws.Range("J136").Value = 2644.25
$ws.Range("K136").Value = 7499.045999999999
$ws.Range("L136").Value = 7932.75
$ws.Range("M136").Value = -4949.045999999999
$ws.Range("N136").Value = -13032.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 9599.416999999999
$ws.Range("I107").Value = 9111.5
$ws.Range("J107").Value = 9843.375
$ws.Range("K107").Value = 9111.5
$ws.Range("L107").Value = 9843.375
$ws.Range("M107").Value = -7191.5
$ws.Range("N107").Value = -13683.375

$ws.Range("H134").Value = 2904.8333
$ws.Range("I134").Value = 2922.8572
$ws.Range("K134").Value = 8768.571599999999
$ws.Range("M134").Value = -6233.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 273.8
$ws.Range("I26").Value = 273.8
$ws.Range("K26").Value = 821.4000000000001
$ws.Range("M26").Value = -533.4000000000001

$ws.Range("H97").Value = 231.13333
$ws.Range("J97").Value = 235.92308
$ws.Range("L97").Value = 707.76924
$ws.Range("N97").Value = -1699.76924

$ws.Range("H113").Value = 1411.125
$ws.Range("I113").Value = 535
$ws.Range("J113").Value = 2287.25
$ws.Range("K113").Value = 1605
$ws.Range("L113").Value = 6861.75
$ws.Range("M113").Value = 565
$ws.Range("N113").Value = -11201.75

$ws.Range("H121").Value = 827.4286
$ws.Range("J121").Value = 926.5
$ws.Range("L121").Value = 2779.5
$ws.Range("N121").Value = -5399.5

$ws.Range("H130").Value = 2642.3333
$ws.Range("I130").Value = 2642.3333
$ws.Range("K130").Value = 7926.999899999999
$ws.Range("M130").Value = -2906.999899999999

$ws.Range("H131").Value = 3615506
$ws.Range("J131").Value = 5130182
$ws.Range("L131").Value = 15390546
$ws.Range("N131").Value = -15400626

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 70000
$ws.Range("J74").Value = 70000
$ws.Range("L74").Value = 70000
$ws.Range("N74").Value = -71872

$ws.Range("H77").Value = 70000
$ws.Range("J77").Value = 70000
$ws.Range("L77").Value = 210000
$ws.Range("N77").Value = -219360

$ws.Range("H132").Value = 2477.4
$ws.Range("I132").Value = 2346.7222
$ws.Range("J132").Value = 2813.4285
$ws.Range("K132").Value = 7040.1666
$ws.Range("L132").Value = 8440.2855
$ws.Range("M132").Value = -4510.1666
$ws.Range("N132").Value = -13500.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 84825
$ws.Range("I141").Value = 84650
$ws.Range("K141").Value = 84650
$ws.Range("M141").Value = -79470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 180091.2
$ws.Range("J64").Value = 180091.2
$ws.Range("L64").Value = 180091.2
$ws.Range("N64").Value = -180587.2

$ws.Range("H67").Value = 180091.2
$ws.Range("J67").Value = 180091.2
$ws.Range("L67").Value = 180091.2
$ws.Range("N67").Value = -181807.2

$ws.Range("H70").Value = 41502.855
$ws.Range("J70").Value = 41502.855
$ws.Range("L70").Value = 41502.855
$ws.Range("N70").Value = -42132.855

$ws.Range("H73").Value = 41502.855
$ws.Range("J73").Value = 41502.855
$ws.Range("L73").Value = 41502.855
$ws.Range("N73").Value = -43686.855

$ws.Range("H107").Value = 2104.5
$ws.Range("J107").Value = 499.5
$ws.Range("L107").Value = 1498.5
$ws.Range("N107").Value = -5338.5

$ws.Range("H126").Value = 3248.1667
$ws.Range("I126").Value = 3179.818
$ws.Range("K126").Value = 9539.454000000002
$ws.Range("M126").Value = -7069.454000000002

$ws.Range("H132").Value = 1790.6072
$ws.Range("I132").Value = 1600.5217
$ws.Range("J132").Value = 2665
$ws.Range("K132").Value = 4801.5651
$ws.Range("L132").Value = 7995
$ws.Range("M132").Value = -2271.5651
$ws.Range("N132").Value = -13055

$ws.Range("H136").Value = 2448
$ws.Range("I136").Value = 2119.3125
$ws.Range("K136").Value = 6357.9375
$ws.Range("M136").Value = -3807.9375
